$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.07953732872513224
$ws.Range("D2").Value = 0.937324309546592

$ws.Range("C3").Value = -0.09503425902082827
$ws.Range("D3").Value = 0.9251479219192009

$ws.Range("C4").Value = 0.2382568997674029
$ws.Range("D4").Value = 0.8138870158302365

$ws.Range("C5").Value = 0.8747107334703133
$ws.Range("D5").Value = 0.391180787095105

$ws.Range("C6").Value = -0.1795906341622638
$ws.Range("D6").Value = 0.8591177847072888

$ws.Range("C7").Value = 0.144578313101597
$ws.Range("D7").Value = 0.886360131452109

$ws.Range("C8").Value = 0.8476746609530327
$ws.Range("D8").Value = 0.4057472884788678

$ws.Range("C9").Value = 0.4119945905521116
$ws.Range("D9").Value = 0.6843300700265322

$ws.Range("C10").Value = 0.6362538055528635
$ws.Range("D10").Value = 0.5311754207076418

$ws.Range("C11").Value = 0.5464115406181709
$ws.Range("D11").Value = 0.5902822695379255
